$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple "Taxonsorteringsordning" (+1) bumps on rows that are not
# --- otherwise touched by the two row swaps below. ---
$ws.Range("B3").Value  = 79245
$ws.Range("B4").Value  = 80350
$ws.Range("B7").Value  = 79245
$ws.Range("B8").Value  = 79245
$ws.Range("B9").Value  = 79245
$ws.Range("B10").Value = 79245
$ws.Range("B14").Value = 79245
$ws.Range("B18").Value = 79245
$ws.Range("B22").Value = 79245

# --- Row 11 / Row 13 effectively swap their species data (the source
# --- export re-ordered these two records); row 11's resulting B value
# --- also picks up the same +1 bump described above. ---
$ws.Range("A11").Value = 131167669
$ws.Range("B11").Value = 79245
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("J11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("Q11").Value = 613256
$ws.Range("R11").Value = 6997380
$ws.Range("AC11").Value = ""
$ws.Range("AF11").Value = ""
$ws.Range("AJ11").Value = "tall"
$ws.Range("AK11").Value = "Pinus sylvestris"
$ws.Range("AO11").Value = "Pinus sylvestris"

$ws.Range("A13").Value = 131167651
$ws.Range("B13").Value = 57884
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = "Tretåig hackspett"
$ws.Range("G13").Value = "Picoides tridactylus"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("J13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = "äldre spår"
$ws.Range("Q13").Value = 613285
$ws.Range("R13").Value = 6997537
$ws.Range("AC13").Value = "Äldre ringhack på tall"
$ws.Range("AF13").Value = ""
$ws.Range("AJ13").Value = ""
$ws.Range("AK13").Value = ""
$ws.Range("AO13").Value = ""

# --- Row 16 / Row 17 likewise swap their species data. ---
$ws.Range("A16").Value = 131167658
$ws.Range("B16").Value = 57073
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 100138
$ws.Range("F16").Value = "Tjäder"
$ws.Range("G16").Value = "Tetrao urogallus"
$ws.Range("H16").Value = "Linnaeus, 1758"
$ws.Range("M16").Value = "färsk spillning"
$ws.Range("Q16").Value = 613330
$ws.Range("R16").Value = 6997326
$ws.Range("AC16").Value = ""

$ws.Range("A17").Value = 131167655
$ws.Range("B17").Value = 57884
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 100109
$ws.Range("F17").Value = "Tretåig hackspett"
$ws.Range("G17").Value = "Picoides tridactylus"
$ws.Range("H17").Value = "(Linnaeus, 1758)"
$ws.Range("M17").Value = "färska spår"
$ws.Range("Q17").Value = 613285
$ws.Range("R17").Value = 6997398
$ws.Range("AC17").Value = "Färska ringhack på tall"
